$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = 45991
$ws.Range("A42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B42").Value = 93
$ws.Range("C42").Value = 106
$ws.Range("D42").Value = 100
